# Weekly crime-data refresh for the CompStat 014 Precinct worksheet.
# Update volume/date header text and refresh all count / % change figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (in-place edits of the rich-text runs) ---
# A8: "Volume 31   Number  34" -> "Volume 31   Number  35"
$ws.Range("A8").Characters(21, 2).Text = "35"

# C9: "Report Covering the Week  8/19/2024  Through  8/25/2024"
#  -> "Report Covering the Week  8/26/2024  Through  9/1/2024"
# (Replace the later date first so the earlier offset stays valid.)
$ws.Range("C9").Characters(47, 9).Text = "9/1/2024"
$ws.Range("C9").Characters(27, 9).Text = "8/26/2024"

# --- Cells switching from a text placeholder ("0"/"***.*") to a real number ---
# D27 and E27 had no prior-year data (shown as "0" / "***.*"); now they do.
$ws.Range("C16").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("E16").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100

# --- Cells switching from a real number back to a text placeholder ---
# D31 & F31 -> "0"; E31 -> "***.*"; C33 -> "0"
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))
$ws.Range("C14").Copy($ws.Range("F31"))
$ws.Range("C14").Copy($ws.Range("C33"))

# --- Remaining numeric updates across the data rows ---
$ws.Range("L14").Value = -33.333333333333
$ws.Range("N14").Value = -55.555555555555
$ws.Range("N15").Value = -29.629629629629
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -55.555555555555
$ws.Range("F16").Value = 35
$ws.Range("G16").Value = 43
$ws.Range("H16").Value = -18.60465116279
$ws.Range("I16").Value = 274
$ws.Range("J16").Value = 324
$ws.Range("K16").Value = -15.432098765432
$ws.Range("L16").Value = -33.170731707317
$ws.Range("M16").Value = 151.376146788991
$ws.Range("N16").Value = -84.431818181818
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -46.153846153846
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = -17.777777777777
$ws.Range("I17").Value = 366
$ws.Range("J17").Value = 342
$ws.Range("K17").Value = 7.017543859649
$ws.Range("L17").Value = 24.067796610169
$ws.Range("M17").Value = 177.272727272727
$ws.Range("N17").Value = -18.666666666666
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 250
$ws.Range("F18").Value = 34
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 252
$ws.Range("J18").Value = 278
$ws.Range("K18").Value = -9.352517985611
$ws.Range("L18").Value = -44.615384615384
$ws.Range("M18").Value = 15.596330275229
$ws.Range("N18").Value = -85.866517106001
$ws.Range("C19").Value = 37
$ws.Range("D19").Value = 42
$ws.Range("E19").Value = -11.904761904761
$ws.Range("F19").Value = 152
$ws.Range("G19").Value = 157
$ws.Range("H19").Value = -3.184713375796
$ws.Range("I19").Value = 1357
$ws.Range("J19").Value = 1575
$ws.Range("K19").Value = -13.841269841269
$ws.Range("L19").Value = -8.31081081081
$ws.Range("M19").Value = -10.310641110376
$ws.Range("N19").Value = -78.750391481365
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 40
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = -9.090909090909
$ws.Range("L20").Value = -16.666666666666
$ws.Range("M20").Value = 166.666666666667
$ws.Range("N20").Value = -82.905982905982
$ws.Range("C21").Value = 57
$ws.Range("D21").Value = 67
$ws.Range("E21").Value = -14.925373134328
$ws.Range("F21").Value = 267
$ws.Range("G21").Value = 269
$ws.Range("H21").Value = -0.743494423791
$ws.Range("I21").Value = 2312
$ws.Range("J21").Value = 2575
$ws.Range("K21").Value = -10.213592233009
$ws.Range("L21").Value = -14.654854189737
$ws.Range("M21").Value = 15.889724310776
$ws.Range("N21").Value = -78.289041224528
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 12
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = -14.285714285714
$ws.Range("I22").Value = 114
$ws.Range("J22").Value = 145
$ws.Range("K22").Value = -21.379310344827
$ws.Range("L22").Value = -5.785123966942
$ws.Range("M22").Value = 23.91304347826
$ws.Range("C24").Value = 93
$ws.Range("D24").Value = 88
$ws.Range("E24").Value = 5.681818181818
$ws.Range("F24").Value = 338
$ws.Range("G24").Value = 352
$ws.Range("H24").Value = -3.977272727272
$ws.Range("I24").Value = 2936
$ws.Range("J24").Value = 2733
$ws.Range("K24").Value = 7.427735089645
$ws.Range("L24").Value = 34.370709382151
$ws.Range("M24").Value = -6.347687400318
$ws.Range("C25").Value = 75
$ws.Range("D25").Value = 78
$ws.Range("E25").Value = -3.846153846153
$ws.Range("F25").Value = 292
$ws.Range("G25").Value = 304
$ws.Range("H25").Value = -3.947368421052
$ws.Range("I25").Value = 2591
$ws.Range("J25").Value = 2501
$ws.Range("K25").Value = 3.598560575769
$ws.Range("L25").Value = 28.07711319822
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 80
$ws.Range("G26").Value = 86
$ws.Range("H26").Value = -6.976744186046
$ws.Range("I26").Value = 700
$ws.Range("J26").Value = 724
$ws.Range("K26").Value = -3.314917127071
$ws.Range("L26").Value = 19.047619047619
$ws.Range("M26").Value = 79.948586118251
$ws.Range("G27").Value = 2
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 26.315789473684
$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 40
$ws.Range("F28").Value = 21
$ws.Range("G28").Value = 20
$ws.Range("H28").Value = 5
$ws.Range("I28").Value = 154
$ws.Range("J28").Value = 154
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 1.315789473684
$ws.Range("H31").Value = -100
